{"js": "// Office.js (Word JavaScript API) edit script.\n// Goal (per commit \"Google analytics second certificate\"):\n//   In the paragraph beginning \"My aim from Coursera is to finish this\n//   specialization...\", change:\n//     \"...I am also in the journey of the first course of the\n//      specialization i.e. Foundations: data, data everywhere and I am\n//      having fun...\"\n//   to:\n//     \"...I am also in the journey of the second course of the\n//      specialization i.e. Ask Questions to Make Data-Driven Decisions\n//      and I am having fun...\"\n\nconst body = context.document.body;\n\n// 1) \"first course of the specialization\" -> \"second course of the specialization\"\nconst firstHits = body.search(\"first course of the specialization\", { matchCase: true });\nfirstHits.load(\"items/text\");\nawait context.sync();\n\nfor (const hit of firstHits.items) {\n  hit.insertText(\"second course of the specialization\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"Foundations: data, data everywhere\" -> \"Ask Questions to Make Data-Driven Decisions\"\nconst secondHits = body.search(\"Foundations: data, data everywhere\", { matchCase: true });\nsecondHits.load(\"items/text\");\nawait context.sync();\n\nfor (const hit of secondHits.items) {\n  hit.insertText(\"Ask Questions to Make Data-Driven Decisions\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Goal (per commit \"Google analytics second certificate\"):\n#   In the paragraph beginning \"My aim from Coursera is to finish this\n#   specialization...\", change:\n#     \"...I am also in the journey of the first course of the\n#      specialization i.e. Foundations: data, data everywhere and I am\n#      having fun...\"\n#   to:\n#     \"...I am also in the journey of the second course of the\n#      specialization i.e. Ask Questions to Make Data-Driven Decisions\n#      and I am having fun...\"\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"first course of the specialization\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"second course of the specialization\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Foundations: data, data everywhere\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Ask Questions to Make Data-Driven Decisions\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
